$wb = $excel.ActiveWorkbook

# Update the "Indrajeet Singh" entry on the Users sheet to "James Craven"
$usersWs = $wb.Worksheets.Item("Users")
$usersWs.Range("A2").Value = "James Craven"

# Move the active sheet/selection from "SaveActivityPopUpMsg" to "Users",
# selecting cell C3 there.
$usersWs.Activate()
$usersWs.Range("C3").Select()
